# Auto-generated edit script: update precision/recall/F1/SEM metrics
# in the "Features" sheet (rows 2-69) and the "Global Metrics" sheet (row 2).
# Values are stored as literal text (comma-decimal locale strings), matching
# the original inline-string cell encoding, so each write forces Text format
# first (NumberFormat "@") and then clears the format residue so the cell
# keeps the workbook default style (no explicit s= style), just like before.

$wb = $excel.ActiveWorkbook
$wsFeatures = $wb.Worksheets.Item("Features")
$wsGlobal = $wb.Worksheets.Item("Global Metrics")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Features sheet ---
Set-TextValue $wsFeatures.Range("B2") "0,707"
Set-TextValue $wsFeatures.Range("C2") "0,824"
Set-TextValue $wsFeatures.Range("D2") "0,761"
Set-TextValue $wsFeatures.Range("E2") "0,997"
Set-TextValue $wsFeatures.Range("B3") "0,737"
Set-TextValue $wsFeatures.Range("C3") "0,890"
Set-TextValue $wsFeatures.Range("D3") "0,807"
Set-TextValue $wsFeatures.Range("E3") "0,985"
Set-TextValue $wsFeatures.Range("B4") "0,719"
Set-TextValue $wsFeatures.Range("C4") "0,863"
Set-TextValue $wsFeatures.Range("D4") "0,784"
Set-TextValue $wsFeatures.Range("E4") "0,997"
Set-TextValue $wsFeatures.Range("B5") "0,729"
Set-TextValue $wsFeatures.Range("C5") "0,909"
Set-TextValue $wsFeatures.Range("D5") "0,809"
Set-TextValue $wsFeatures.Range("E5") "1,000"
Set-TextValue $wsFeatures.Range("B6") "0,468"
Set-TextValue $wsFeatures.Range("C6") "0,463"
Set-TextValue $wsFeatures.Range("D6") "0,466"
Set-TextValue $wsFeatures.Range("E6") "1,000"
Set-TextValue $wsFeatures.Range("B7") "0,446"
Set-TextValue $wsFeatures.Range("C7") "0,466"
Set-TextValue $wsFeatures.Range("D7") "0,456"
Set-TextValue $wsFeatures.Range("E7") "1,000"
Set-TextValue $wsFeatures.Range("B8") "0,328"
Set-TextValue $wsFeatures.Range("C8") "0,352"
Set-TextValue $wsFeatures.Range("D8") "0,339"
Set-TextValue $wsFeatures.Range("E8") "0,966"
Set-TextValue $wsFeatures.Range("B9") "0,517"
Set-TextValue $wsFeatures.Range("C9") "0,756"
Set-TextValue $wsFeatures.Range("D9") "0,614"
Set-TextValue $wsFeatures.Range("E9") "1,000"
Set-TextValue $wsFeatures.Range("B10") "0,611"
Set-TextValue $wsFeatures.Range("C10") "0,868"
Set-TextValue $wsFeatures.Range("D10") "0,717"
Set-TextValue $wsFeatures.Range("E10") "0,974"
Set-TextValue $wsFeatures.Range("B11") "0,568"
Set-TextValue $wsFeatures.Range("C11") "0,553"
Set-TextValue $wsFeatures.Range("D11") "0,560"
Set-TextValue $wsFeatures.Range("E11") "0,977"
Set-TextValue $wsFeatures.Range("B12") "0,568"
Set-TextValue $wsFeatures.Range("C12") "0,583"
Set-TextValue $wsFeatures.Range("D12") "0,575"
Set-TextValue $wsFeatures.Range("E12") "0,957"
Set-TextValue $wsFeatures.Range("B13") "0,605"
Set-TextValue $wsFeatures.Range("C13") "0,561"
Set-TextValue $wsFeatures.Range("D13") "0,582"
Set-TextValue $wsFeatures.Range("E13") "0,978"
Set-TextValue $wsFeatures.Range("B14") "0,632"
Set-TextValue $wsFeatures.Range("C14") "0,615"
Set-TextValue $wsFeatures.Range("D14") "0,623"
Set-TextValue $wsFeatures.Range("E14") "0,950"
Set-TextValue $wsFeatures.Range("B15") "0,667"
Set-TextValue $wsFeatures.Range("C15") "0,821"
Set-TextValue $wsFeatures.Range("D15") "0,736"
Set-TextValue $wsFeatures.Range("E15") "1,000"
Set-TextValue $wsFeatures.Range("B16") "0,453"
Set-TextValue $wsFeatures.Range("C16") "0,706"
Set-TextValue $wsFeatures.Range("D16") "0,552"
Set-TextValue $wsFeatures.Range("E16") "1,000"
Set-TextValue $wsFeatures.Range("B17") "0,451"
Set-TextValue $wsFeatures.Range("C17") "0,697"
Set-TextValue $wsFeatures.Range("D17") "0,548"
Set-TextValue $wsFeatures.Range("E17") "1,000"
Set-TextValue $wsFeatures.Range("B18") "0,519"
Set-TextValue $wsFeatures.Range("C18") "0,583"
Set-TextValue $wsFeatures.Range("D18") "0,549"
Set-TextValue $wsFeatures.Range("E18") "0,959"
Set-TextValue $wsFeatures.Range("B19") "0,600"
Set-TextValue $wsFeatures.Range("C19") "0,720"
Set-TextValue $wsFeatures.Range("D19") "0,655"
Set-TextValue $wsFeatures.Range("E19") "0,850"
Set-TextValue $wsFeatures.Range("B20") "0,600"
Set-TextValue $wsFeatures.Range("C20") "0,750"
Set-TextValue $wsFeatures.Range("D20") "0,667"
Set-TextValue $wsFeatures.Range("E20") "0,819"
Set-TextValue $wsFeatures.Range("B21") "0,524"
Set-TextValue $wsFeatures.Range("C21") "0,524"
Set-TextValue $wsFeatures.Range("D21") "0,524"
Set-TextValue $wsFeatures.Range("E21") "0,845"
Set-TextValue $wsFeatures.Range("B22") "0,857"
Set-TextValue $wsFeatures.Range("C22") "0,375"
Set-TextValue $wsFeatures.Range("D22") "0,522"
Set-TextValue $wsFeatures.Range("E22") "1,000"
Set-TextValue $wsFeatures.Range("B23") "0,400"
Set-TextValue $wsFeatures.Range("C23") "0,471"
Set-TextValue $wsFeatures.Range("D23") "0,432"
Set-TextValue $wsFeatures.Range("E23") "0,898"
Set-TextValue $wsFeatures.Range("B24") "0,857"
Set-TextValue $wsFeatures.Range("C24") "0,429"
Set-TextValue $wsFeatures.Range("D24") "0,571"
Set-TextValue $wsFeatures.Range("E24") "0,940"
Set-TextValue $wsFeatures.Range("B25") "0,500"
Set-TextValue $wsFeatures.Range("C25") "0,588"
Set-TextValue $wsFeatures.Range("D25") "0,541"
Set-TextValue $wsFeatures.Range("E25") "0,892"
Set-TextValue $wsFeatures.Range("B26") "0,400"
Set-TextValue $wsFeatures.Range("C26") "0,308"
Set-TextValue $wsFeatures.Range("D26") "0,348"
Set-TextValue $wsFeatures.Range("E26") "1,000"
Set-TextValue $wsFeatures.Range("B27") "0,579"
Set-TextValue $wsFeatures.Range("C27") "0,524"
Set-TextValue $wsFeatures.Range("D27") "0,550"
Set-TextValue $wsFeatures.Range("E27") "1,000"
Set-TextValue $wsFeatures.Range("B28") "0,297"
Set-TextValue $wsFeatures.Range("C28") "0,688"
Set-TextValue $wsFeatures.Range("D28") "0,415"
Set-TextValue $wsFeatures.Range("E28") "1,000"
Set-TextValue $wsFeatures.Range("B29") "0,167"
Set-TextValue $wsFeatures.Range("C29") "0,333"
Set-TextValue $wsFeatures.Range("D29") "0,222"
Set-TextValue $wsFeatures.Range("E29") "1,000"
Set-TextValue $wsFeatures.Range("B30") "0,455"
Set-TextValue $wsFeatures.Range("C30") "0,417"
Set-TextValue $wsFeatures.Range("D30") "0,435"
Set-TextValue $wsFeatures.Range("E30") "1,000"
Set-TextValue $wsFeatures.Range("B31") "0,188"
Set-TextValue $wsFeatures.Range("C31") "0,300"
Set-TextValue $wsFeatures.Range("D31") "0,231"
Set-TextValue $wsFeatures.Range("E31") "0,935"
Set-TextValue $wsFeatures.Range("B32") "0,207"
Set-TextValue $wsFeatures.Range("C32") "0,545"
Set-TextValue $wsFeatures.Range("D32") "0,300"
Set-TextValue $wsFeatures.Range("E32") "1,000"
Set-TextValue $wsFeatures.Range("B33") "0,692"
Set-TextValue $wsFeatures.Range("C33") "0,818"
Set-TextValue $wsFeatures.Range("D33") "0,750"
Set-TextValue $wsFeatures.Range("E33") "1,000"
Set-TextValue $wsFeatures.Range("B34") "0,556"
Set-TextValue $wsFeatures.Range("C34") "0,417"
Set-TextValue $wsFeatures.Range("D34") "0,476"
Set-TextValue $wsFeatures.Range("E34") "1,000"
Set-TextValue $wsFeatures.Range("B35") "0,357"
Set-TextValue $wsFeatures.Range("C35") "0,500"
Set-TextValue $wsFeatures.Range("D35") "0,417"
Set-TextValue $wsFeatures.Range("E35") "1,000"
Set-TextValue $wsFeatures.Range("B36") "0,750"
Set-TextValue $wsFeatures.Range("C36") "0,429"
Set-TextValue $wsFeatures.Range("D36") "0,545"
Set-TextValue $wsFeatures.Range("E36") "1,000"
Set-TextValue $wsFeatures.Range("B37") "0,222"
Set-TextValue $wsFeatures.Range("C37") "0,222"
Set-TextValue $wsFeatures.Range("D37") "0,222"
Set-TextValue $wsFeatures.Range("E37") "1,000"
Set-TextValue $wsFeatures.Range("B38") "0,300"
Set-TextValue $wsFeatures.Range("C38") "0,429"
Set-TextValue $wsFeatures.Range("D38") "0,353"
Set-TextValue $wsFeatures.Range("E38") "1,000"
Set-TextValue $wsFeatures.Range("B39") "0,667"
Set-TextValue $wsFeatures.Range("C39") "0,500"
Set-TextValue $wsFeatures.Range("D39") "0,571"
Set-TextValue $wsFeatures.Range("E39") "1,000"
Set-TextValue $wsFeatures.Range("B40") "0,600"
Set-TextValue $wsFeatures.Range("C40") "0,750"
Set-TextValue $wsFeatures.Range("D40") "0,667"
Set-TextValue $wsFeatures.Range("E40") "1,000"
Set-TextValue $wsFeatures.Range("B41") "0,625"
Set-TextValue $wsFeatures.Range("C41") "0,714"
Set-TextValue $wsFeatures.Range("D41") "0,667"
Set-TextValue $wsFeatures.Range("E41") "1,000"
Set-TextValue $wsFeatures.Range("B43") "0,500"
Set-TextValue $wsFeatures.Range("C43") "0,667"
Set-TextValue $wsFeatures.Range("D43") "0,571"
Set-TextValue $wsFeatures.Range("E43") "1,000"
Set-TextValue $wsFeatures.Range("B44") "0,500"
Set-TextValue $wsFeatures.Range("C44") "0,333"
Set-TextValue $wsFeatures.Range("D44") "0,400"
Set-TextValue $wsFeatures.Range("E44") "0,897"
Set-TextValue $wsFeatures.Range("E45") "1,000"
Set-TextValue $wsFeatures.Range("B46") "0,333"
Set-TextValue $wsFeatures.Range("C46") "0,500"
Set-TextValue $wsFeatures.Range("D46") "0,400"
Set-TextValue $wsFeatures.Range("E46") "1,000"
Set-TextValue $wsFeatures.Range("B47") "0,625"
Set-TextValue $wsFeatures.Range("C47") "0,556"
Set-TextValue $wsFeatures.Range("D47") "0,588"
Set-TextValue $wsFeatures.Range("E47") "1,000"
Set-TextValue $wsFeatures.Range("B48") "0,625"
Set-TextValue $wsFeatures.Range("C48") "0,625"
Set-TextValue $wsFeatures.Range("D48") "0,625"
Set-TextValue $wsFeatures.Range("E48") "1,000"
Set-TextValue $wsFeatures.Range("B49") "0,333"
Set-TextValue $wsFeatures.Range("C49") "0,500"
Set-TextValue $wsFeatures.Range("D49") "0,400"
Set-TextValue $wsFeatures.Range("E49") "1,000"
Set-TextValue $wsFeatures.Range("B50") "1,000"
Set-TextValue $wsFeatures.Range("C50") "0,667"
Set-TextValue $wsFeatures.Range("D50") "0,800"
Set-TextValue $wsFeatures.Range("E50") "0,667"
Set-TextValue $wsFeatures.Range("B53") "0,250"
Set-TextValue $wsFeatures.Range("C53") "0,200"
Set-TextValue $wsFeatures.Range("D53") "0,222"
Set-TextValue $wsFeatures.Range("E53") "0,914"
Set-TextValue $wsFeatures.Range("B55") "0,143"
Set-TextValue $wsFeatures.Range("C55") "0,333"
Set-TextValue $wsFeatures.Range("D55") "0,200"
Set-TextValue $wsFeatures.Range("E55") "1,000"
Set-TextValue $wsFeatures.Range("B56") "1,000"
Set-TextValue $wsFeatures.Range("C56") "1,000"
Set-TextValue $wsFeatures.Range("D56") "1,000"
Set-TextValue $wsFeatures.Range("E56") "1,000"
Set-TextValue $wsFeatures.Range("E60") "0,463"
Set-TextValue $wsFeatures.Range("E61") "0,643"
Set-TextValue $wsFeatures.Range("B69") "1,000"
Set-TextValue $wsFeatures.Range("C69") "0,500"
Set-TextValue $wsFeatures.Range("D69") "0,667"
Set-TextValue $wsFeatures.Range("E69") "0,500"

# --- Global Metrics sheet ---
Set-TextValue $wsGlobal.Range("B2") "0,398"
Set-TextValue $wsGlobal.Range("C2") "0,739"
Set-TextValue $wsGlobal.Range("D2") "0,575"
Set-TextValue $wsGlobal.Range("E2") "0,952"
